$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: 四方坪站 (site 1)
$ws.Range("A10").Value = 45996
$ws.Range("B10").Value = "四方坪站"
$ws.Range("C10").Value = 8708.01
$ws.Range("D10").Value = 7611.33
$ws.Range("E10").Value = 2907.12
$ws.Range("F10").Value = 370

# Row 11: 高岭站 (site 2)
$ws.Range("A11").Value = 45996
$ws.Range("B11").Value = "高岭站"
$ws.Range("C11").Value = 5598.96
$ws.Range("D11").Value = 4780.25
$ws.Range("E11").Value = 1584.5
$ws.Range("F11").Value = 188

$ws.Range("I8").Select()
